$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric to Excel; force text format
# first so the literal string (including trailing/leading zeros) is preserved,
# matching the inline-string cells in the source workbook.
$textRows = @(5,6,8,10,11,14,16,18,19,20,21,22,23,24,26,27,28,29,30,31,33,34,35,36,37,38,39,40,41,42,43,45,46,47,48,49)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "66.855.71"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "3.494.54"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "593.54"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "172.52"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "0.432"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "4.105.17"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "28.85"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "66.944.17"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "3.498.68"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "6.29"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "14.08"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "394.12"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "7.99"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "72.86"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("D24").Value = "0.535"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "10.23"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "6.24"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "23.68"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "7.35"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "162.97"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "0.880"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").Value = "6.93"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "4.66"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").Value = "0.0743"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "27.34"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").Value = "26.22"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "2.800.73"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D47").Value = "0.0302"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "337.05"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").Value = "34.23"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("E51").Value = "  -0.64%  "

# Row 45/46: the coin order swapped (OKB now ranks above dogwifhat),
# along with updated price/volume figures.
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "42.62"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  +1.42%  "
